# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# The "K" column (G) values are recalculated/regenerated; write the new values
# for rows 2-19 in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 0
    6  = 2
    7  = 6
    8  = 12
    9  = 2
    10 = 3
    11 = 6
    12 = 3
    13 = 3
    14 = 1
    15 = 2
    16 = 1
    17 = 1
    18 = 5
    19 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
